# The commit swaps the contents of ppt/theme/theme1.xml (the deck's
# active/master theme, originally the "Integral" / "Red Violet" palette)
# and ppt/theme/theme2.xml (an orphaned theme used only by the notes
# master, originally the "Office Theme" palette) - theme1 ends up with
# the Office Theme colors and theme2 ends up with the Integral colors.
#
# The PowerPoint object model only exposes a single, deck-wide theme
# (reachable from the slide master / any slide's ThemeColorScheme), which
# always maps onto ppt/theme/theme1.xml, so we reproduce the
# user-visible part of that swap here: push the "Office Theme" color
# scheme (the 12 standard theme colors) onto the presentation's theme.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

# MsoThemeColorSchemeIndex order: 1 Dark1, 2 Light1, 3 Dark2, 4 Light2,
# 5 Accent1 .. 10 Accent6, 11 Hyperlink, 12 FollowedHyperlink.
$tcs.Item(1).RGB  = "000000"  # dk1
$tcs.Item(2).RGB  = "FFFFFF"  # lt1
$tcs.Item(3).RGB  = "44546A"  # dk2
$tcs.Item(4).RGB  = "E7E6E6"  # lt2
$tcs.Item(5).RGB  = "5B9BD5"  # accent1
$tcs.Item(6).RGB  = "ED7D31"  # accent2
$tcs.Item(7).RGB  = "A5A5A5"  # accent3
$tcs.Item(8).RGB  = "FFC000"  # accent4
$tcs.Item(9).RGB  = "4472C4"  # accent5
$tcs.Item(10).RGB = "70AD47"  # accent6
$tcs.Item(11).RGB = "0563C1"  # hlink
$tcs.Item(12).RGB = "954F72"  # folHlink
